$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the current title text (currently in B1) before we move it.
$title = $ws.Range("B1").Value2

# New row 1: "Unit:" / "persons"
$ws.Range("A1").Value = "Unit:"

# New row 3: "Source:" / "[1]" / "Eurostat" / <title>
$ws.Range("A3").Value = "Source:"
$ws.Range("B3").Value = "[1]"
$ws.Range("C3").Value = "Eurostat"
$ws.Range("D3").Value = $title

# "persons" goes into B1
$ws.Range("B1").Value = "persons"

# The title is no longer bold now that it lives in D3; clear bold everywhere it used to be.
$ws.Range("B1").Font.Bold = $false
$ws.Range("D3").Font.Bold = $false

# Match the new selection / page setup from the authored workbook.
$ws.Range("B2").Select() | Out-Null
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1 | Out-Null
